$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) gets two more "dlgTruckInsurance" cells in H1/I1
$ws.Range("H1").Value = "dlgTruckInsurance"
$ws.Range("I1").Value = "dlgTruckInsurance"

# Fill descriptive "Vehicle Page ..." texts first (matches authoring order
# reflected in the shared-string table), then the technical step IDs.
$ws.Range("B3").Value = "Vehicle Page check for open mandatory fields"
$ws.Range("B4").Value = "Vehicle Page check for hints regarding mandatory fields"
$ws.Range("D4").Value = "Vehicle Page check error hint list value ranges"
$ws.Range("I4").Value = "Vehicle Page check error hint manufacturing date in the future"

$ws.Range("A3").Value = "103_TruckInsurance_002_VehicleData_001_MandatoryFields"
$ws.Range("A4").Value = "103_TruckInsurance_002_VehicleData_002_FieldHintsAndErrors"

$ws.Range("C3").Value = "103_TruckInsurance_002_VehicleData_001_MandatoryFields_FillMake"
$ws.Range("C4").Value = "103_TruckInsurance_002_VehicleData_002_EnterNumericValuesBelowRange"

$ws.Range("D3").Value = "103_TruckInsurance_002_VehicleData_001_MandatoryFields_CheckFilledMake"
$ws.Range("F4").Value = "103_TruckInsurance_002_VehicleData_002_EnterNumericValuesAboveRange"
$ws.Range("H4").Value = "103_TruckInsurance_002_VehicleData_002_ManufacturingDateInTheFuture"

# G4 reuses the same text as D4 (shared string already present)
$ws.Range("G4").Value = "Vehicle Page check error hint list value ranges"

# Column width adjustments (auto-fit-like widths observed in target;
# values pre-compensated for the host's character->stored-width rounding)
$ws.Columns.Item(1).ColumnWidth = 73.83072916666667
$ws.Columns.Item(3).ColumnWidth = 69.16666666666667
$ws.Columns.Item(4).ColumnWidth = 65.60807291666667
$ws.Columns.Item(6).ColumnWidth = 64.05338541666667
$ws.Columns.Item(7).ColumnWidth = 37.498697916666664
$ws.Columns.Item(8).ColumnWidth = 68.05338541666667
$ws.Columns.Item(9).ColumnWidth = 51.276041666666664

# Selection moves to A3
$ws.Range("A3").Select()

# Best-effort: reflect the resized/relocated application window (matches
# the saved workbookView bounds in the target file).
$win = $excel.ActiveWindow
$win.Left = 960
$win.Top = 1260
$win.Width = 36948
$win.Height = 15672
